$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = -0.011155
$ws.Range("E2").Value = -0.09970000000000001
$ws.Range("G2").Value = 0.4213893239189871
$ws.Range("H2").Value = 0.4213893239189871
$ws.Range("I2").Value = 0.3587449897804794
$ws.Range("J2").Value = 0.3116639388897583
$ws.Range("K2").Value = 133.86
$ws.Range("L2").Value = 0.1776603933851831
$ws.Range("M2").Value = 113.82
$ws.Range("N2").Value = 0.1213303485769108
$ws.Range("O2").Value = 0.8502913491707753
$ws.Range("P2").Value = 113.82
$ws.Range("Q2").Value = 0.1213303485769108
$ws.Range("R2").Value = 0.8502913491707753
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 550.5
$ws.Range("V2").Value = 0.5868244323632875
$ws.Range("W2").Value = 0.04121709732147109
$ws.Range("X2").Value = 0.08182782500622357
$ws.Range("Y2").Value = -0.04061072768475248
$ws.Range("Z2").Value = 0.2165525645668431
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.05022363700749735
$ws.Range("AC2").Value = -0.04952077106789983
$ws.Range("AD2").Value = 2480.6
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 2480.6
$ws.Range("AG2").Value = 1930.1
$ws.Range("AH2").Value = 0.725597449322842
$ws.Range("AI2").Value = 0.6243330313097755
$ws.Range("AJ2").Value = 0.672930757966669
$ws.Range("AK2").Value = 0.5639115318315949
$ws.Range("AL2").Value = 138.8
$ws.Range("AM2").Value = 138.8
$ws.Range("AN2").Value = 8.719156414762741
$ws.Range("AO2").Value = 1.947406340057637
$ws.Range("AP2").Value = 6.784182776801406
$ws.Range("AQ2").Value = 1.947406340057637

# Row 3
$ws.Range("B3").Value = 'Oman International Development and Investment Company SAOG (MSM:OMVS)'
$ws.Range("D3").Value = 0.208
$ws.Range("E3").Value = 0.154
$ws.Range("G3").Value = 0.4692134107027724
$ws.Range("H3").Value = 0.4692134107027724
$ws.Range("I3").Value = 0.4021598968407479
$ws.Range("J3").Value = 0.3492243146614605
$ws.Range("K3").Value = 100.4
$ws.Range("L3").Value = 0.1618310767246938
$ws.Range("M3").Value = 83.8
$ws.Range("N3").Value = 0.1428327935912732
$ws.Range("O3").Value = 0.8346613545816732
$ws.Range("P3").Value = 83.8
$ws.Range("Q3").Value = 0.1428327935912732
$ws.Range("R3").Value = 0.8346613545816732
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 324.6
$ws.Range("V3").Value = 0.5532640190898245
$ws.Range("W3").Value = 0.122873577285522
$ws.Range("X3").Value = 0.03829981032240859
$ws.Range("Y3").Value = 0.08457376696311338
$ws.Range("Z3").Value = 0.6191616766467064
$ws.Range("AA3").Value = 0.2162263121915868
$ws.Range("AB3").Value = 0.0464837209993795
$ws.Range("AC3").Value = 0.1697425911922073
$ws.Range("AD3").Value = 795.2
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 795.2
$ws.Range("AG3").Value = 470.6
$ws.Range("AH3").Value = 0.5754396121282292
$ws.Range("AI3").Value = 0.5209984930878595
$ws.Range("AJ3").Value = 0.4450959992433557
$ws.Range("AK3").Value = 0.3916118831655155
$ws.Range("AL3").Value = 122.2
$ws.Range("AM3").Value = 122.2
$ws.Range("AN3").Value = 3.018982536066819
$ws.Range("AO3").Value = 2.041734860883797
$ws.Range("AP3").Value = 1.78663629460896
$ws.Range("AQ3").Value = 2.041734860883797

# Row 4
$ws.Range("D4").Value = -0.0202
$ws.Range("E4").Value = -0.28
$ws.Range("G4").Value = 0.5986394557823129
$ws.Range("H4").Value = 0.5986394557823129
$ws.Range("I4").Value = 0.4716553287981859
$ws.Range("J4").Value = 0.381254724111867
$ws.Range("K4").Value = 3.39
$ws.Range("L4").Value = 0.07687074829931972
$ws.Range("M4").Value = 5.35
$ws.Range("N4").Value = 0.1036821705426356
$ws.Range("O4").Value = 1.578171091445428
$ws.Range("P4").Value = 5.35
$ws.Range("Q4").Value = 0.1036821705426356
$ws.Range("R4").Value = 1.578171091445428
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 191.8
$ws.Range("V4").Value = 3.717054263565891
$ws.Range("W4").Value = 0.02009484291641968
$ws.Range("X4").Value = 0.09167693776384939
$ws.Range("Y4").Value = -0.07158209484742971
$ws.Range("Z4").Value = 0.114456267843239
$ws.Range("AA4").Value = 0.04363699281944805
$ws.Range("AB4").Value = 0.05024615830075058
$ws.Range("AC4").Value = -0.006609165481302526
$ws.Range("AD4").Value = 310.4
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 310.4
$ws.Range("AG4").Value = 118.6
$ws.Range("AH4").Value = 0.8574585635359115
$ws.Range("AI4").Value = 0.6548523206751055
$ws.Range("AJ4").Value = 0.6968272620446533
$ws.Range("AK4").Value = 0.4202693125442948
$ws.Range("AL4").Value = 16.6
$ws.Range("AM4").Value = 16.6
$ws.Range("AN4").Value = 14.71090047393365
$ws.Range("AO4").Value = 1.253012048192771
$ws.Range("AP4").Value = 5.620853080568718
$ws.Range("AQ4").Value = 1.253012048192771

# Row 5
$ws.Range("B5").Value = 'National Finance Company SAOG (MSM:NFCI)'
$ws.Range("D5").Value = 0.152
$ws.Range("E5").Value = 0.129
$ws.Range("K5").Value = 27.5
$ws.Range("L5").Value = 0.5149812734082397
$ws.Range("M5").Value = 19.4
$ws.Range("N5").Value = 0.1104154809334092
$ws.Range("O5").Value = 0.7054545454545454
$ws.Range("P5").Value = 19.4
$ws.Range("Q5").Value = 0.1104154809334092
$ws.Range("R5").Value = 0.7054545454545454
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 17.4
$ws.Range("V5").Value = 0.09903244166192374
$ws.Range("W5").Value = 0.1025736665423349
$ws.Range("X5").Value = 0.07276354469406951
$ws.Range("Y5").Value = 0.02981012184826543
$ws.Range("Z5").Value = 0.04951780415430267
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.04884042642155553
$ws.Range("AC5").Value = -0.04884042642155553
$ws.Range("AD5").Value = 766.8
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 766.8
$ws.Range("AG5").Value = 749.4
$ws.Range("AH5").Value = 0.8135809018567639
$ws.Range("AI5").Value = 0.7353279631760644
$ws.Range("AJ5").Value = 0.8100745865311859
$ws.Range("AK5").Value = 0.7308367466354593

# Row 6
$ws.Range("B6").Value = 'Taageer Finance Company SAOG (MSM:TFCI)'
$ws.Range("D6").Value = -0.00211
$ws.Range("E6").Value = -0.09970000000000001
$ws.Range("K6").Value = 6.35
$ws.Range("L6").Value = 0.3207070707070707
$ws.Range("M6").Value = 5.27
$ws.Range("N6").Value = 0.09634369287020109
$ws.Range("O6").Value = 0.8299212598425196
$ws.Range("P6").Value = 5.27
$ws.Range("Q6").Value = 0.09634369287020109
$ws.Range("R6").Value = 0.8299212598425196
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 6.07
$ws.Range("V6").Value = 0.1109689213893967
$ws.Range("W6").Value = 0.06001890359168242
$ws.Range("X6").Value = 0.09089210531837762
$ws.Range("Y6").Value = -0.0308732017266952
$ws.Range("Z6").Value = 0.04119678748283467
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.05020111571424413
$ws.Range("AC6").Value = -0.05020111571424413
$ws.Range("AD6").Value = 325.3
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 325.3
$ws.Range("AG6").Value = 319.23
$ws.Range("AH6").Value = 0.8560526315789474
$ws.Range("AI6").Value = 0.7528349919000231
$ws.Range("AJ6").Value = 0.8537159361377799
$ws.Range("AK6").Value = 0.7493134286317865

# Row 7
$ws.Range("D7").Value = -0.248
$ws.Range("E7").ClearContents()
$ws.Range("K7").Value = -6.36
$ws.Range("L7").Value = -1.232558139534884
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").ClearContents()
$ws.Range("U7").Value = 3.34
$ws.Range("V7").Value = 0.1132203389830508
$ws.Range("W7").Value = -0.06133076181292189
$ws.Range("X7").Value = 0.0944890466291567
$ws.Range("Y7").Value = -0.1558198084420786
$ws.Range("Z7").Value = 0.0170516506394369
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.05040056897588463
$ws.Range("AC7").Value = -0.05040056897588463
$ws.Range("AD7").Value = 184.7
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 184.7
$ws.Range("AG7").Value = 181.36
$ws.Range("AH7").Value = 0.8622782446311859
$ws.Range("AI7").Value = 0.6545003543586109
$ws.Range("AJ7").Value = 0.8600967466565493
$ws.Range("AK7").Value = 0.6503621889119988

# Row 8
$ws.Range("B8").Value = 'United Finance Company SAOG (MSM:UFCI)'
$ws.Range("D8").Value = -0.147
$ws.Range("E8").Value = -0.283
$ws.Range("K8").Value = 2.58
$ws.Range("L8").Value = 0.2433962264150943
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("T8").ClearContents()
$ws.Range("U8").Value = 7.29
$ws.Range("V8").Value = 0.1827067669172932
$ws.Range("W8").Value = 0.02241529105125977
$ws.Range("X8").Value = 0.05096537071624238
$ws.Range("Y8").Value = -0.02855007966498261
$ws.Range("Z8").Value = 0.04600494770192266
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.05207213264336406
$ws.Range("AC8").Value = -0.05207213264336406
$ws.Range("AD8").Value = 98.2
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 98.2
$ws.Range("AG8").Value = 90.91
$ws.Range("AH8").Value = 0.7110789283128168
$ws.Range("AI8").Value = 0.4550509731232623
$ws.Range("AJ8").Value = 0.6949774482073235
$ws.Range("AK8").Value = 0.4359982734641025
